# Add transformation funs with accessories
#
# - Sheet21: "select all" selection widened from row 1 to the full sheet
#   (A1:XFD1 -> A1:XFD1048576).
# - Sheet22: loses tabSelected (a later sheet becomes the active tab).
# - Two new sheets (Sheet23, Sheet24) are appended, each a header row plus
#   two data rows, cloned from the existing "QOQ" transformation row with
#   small tweaks (an indeks_obdobje year on one, a lower-cased "qoq" code
#   and a "--M" series id on the other).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Fix up the stale "row 1 only" select-all selection on Sheet21.
# ---------------------------------------------------------------------
$ws21 = $wb.Worksheets.Item(21)
$ws21.Activate()
$ws21.Range("A1:XFD1048576").Select()

# ---------------------------------------------------------------------
# 2. Append two new worksheets at the end of the workbook (after the
#    current last sheet, "Sheet22"), matching Excel's auto-naming.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws23 = $wb.Worksheets.Add($null, $lastSheet)

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws24 = $wb.Worksheets.Add($null, $lastSheet2)

# ---------------------------------------------------------------------
# 3. Sheet23 — header row + the "QOQ" series row (now carrying an
#    indeks_obdobje of 2015) + the companion "število" series row.
# ---------------------------------------------------------------------
$ws23.Range("A1").Value = "serija"
$ws23.Range("B1").Value = "enota"
$ws23.Range("C1").Value = "legenda"
$ws23.Range("D1").Value = "barva"
$ws23.Range("E1").Value = "tip"
$ws23.Range("F1").Value = "stacked"
$ws23.Range("G1").Value = "drseca_obdobja"
$ws23.Range("H1").Value = "drseca_poravnava"
$ws23.Range("I1").Value = "rast"
$ws23.Range("J1").Value = "indeks_obdobje"
$ws23.Range("K1").Value = "velikost"
$ws23.Range("L1").Value = "naslov"
$ws23.Range("M1").Value = "xmin"
$ws23.Range("N1").Value = "xmax"
$ws23.Range("O1").Value = "opomba"
$ws23.Range("P1").Value = "stolpci_legende"
$ws23.Range("Q1").Value = "datum_podatkov"
$ws23.Range("R1").Value = "leva_y_os"
$ws23.Range("S1").Value = "desna_y_os"

$ws23.Range("A2").Value = "SURS--0300230S--P31_S14_D--G4--N--Q"
$ws23.Range("B2").Value = "%"
$ws23.Range("C2").Value = "Serija 1"
$ws23.Range("E2").Value = "line"
$ws23.Range("G2").Value = 3
$ws23.Range("H2").Value = "c"
$ws23.Range("I2").Value = "QOQ"
$ws23.Range("J2").Value = 2015
$ws23.Range("K2").Value = 1
$ws23.Range("L2").Value = "Naslov"
$ws23.Range("M2").Value = 36557
$ws23.Range("M2").NumberFormat = "m/d/yy"
$ws23.Range("N2").Value = 12
$ws23.Range("O2").Value = "* nekaj nekaj"
$ws23.Range("P2").Value = 3
$ws23.Range("R2").Value = "1,2,3"
$ws23.Range("S2").Value = "10,20,30"

$ws23.Range("A3").Value = "SURS--0300230S--P31_S15_D--G4--N--Q"
$ws23.Range("B3").Value = "število"
$ws23.Range("C3").Value = "Serija 2"
$ws23.Range("D3").Value = 1
$ws23.Range("E3").Value = "bar"
$ws23.Range("G3").Value = 3
$ws23.Range("H3").Value = "d"
$ws23.Range("K3").Value = 1
$ws23.Range("L3").Value = "Naslov"
$ws23.Range("M3").Value = 36557
$ws23.Range("M3").NumberFormat = "m/d/yy"
$ws23.Range("N3").Value = 12
$ws23.Range("O3").Value = "* nekaj nekaj"
$ws23.Range("P3").Value = 3

$ws23.Range("J2").Select()
$ws23.Range("A1:S3").Select()

# ---------------------------------------------------------------------
# 4. Sheet24 — same shape, but the "rast" code is lower-cased ("qoq")
#    and the series id targets the monthly ("--M") publication instead
#    of the quarterly one, with no indeks_obdobje set.
# ---------------------------------------------------------------------
$ws24.Range("A1").Value = "serija"
$ws24.Range("B1").Value = "enota"
$ws24.Range("C1").Value = "legenda"
$ws24.Range("D1").Value = "barva"
$ws24.Range("E1").Value = "tip"
$ws24.Range("F1").Value = "stacked"
$ws24.Range("G1").Value = "drseca_obdobja"
$ws24.Range("H1").Value = "drseca_poravnava"
$ws24.Range("I1").Value = "rast"
$ws24.Range("J1").Value = "indeks_obdobje"
$ws24.Range("K1").Value = "velikost"
$ws24.Range("L1").Value = "naslov"
$ws24.Range("M1").Value = "xmin"
$ws24.Range("N1").Value = "xmax"
$ws24.Range("O1").Value = "opomba"
$ws24.Range("P1").Value = "stolpci_legende"
$ws24.Range("Q1").Value = "datum_podatkov"
$ws24.Range("R1").Value = "leva_y_os"
$ws24.Range("S1").Value = "desna_y_os"

$ws24.Range("I2").Value = "qoq"
$ws24.Range("A2").Value = "SURS--0300230S--P31_S14_D--G4--N--M"
$ws24.Range("B2").Value = "%"
$ws24.Range("C2").Value = "Serija 1"
$ws24.Range("E2").Value = "line"
$ws24.Range("G2").Value = 3
$ws24.Range("H2").Value = "c"
$ws24.Range("K2").Value = 1
$ws24.Range("L2").Value = "Naslov"
$ws24.Range("M2").Value = 36557
$ws24.Range("M2").NumberFormat = "m/d/yy"
$ws24.Range("N2").Value = 12
$ws24.Range("O2").Value = "* nekaj nekaj"
$ws24.Range("P2").Value = 3
$ws24.Range("R2").Value = "1,2,3"
$ws24.Range("S2").Value = "10,20,30"

$ws24.Range("A3").Value = "SURS--0300230S--P31_S15_D--G4--N--Q"
$ws24.Range("B3").Value = "število"
$ws24.Range("C3").Value = "Serija 2"
$ws24.Range("D3").Value = 1
$ws24.Range("E3").Value = "bar"
$ws24.Range("G3").Value = 3
$ws24.Range("H3").Value = "d"
$ws24.Range("K3").Value = 1
$ws24.Range("L3").Value = "Naslov"
$ws24.Range("M3").Value = 36557
$ws24.Range("M3").NumberFormat = "m/d/yy"
$ws24.Range("N3").Value = 12
$ws24.Range("O3").Value = "* nekaj nekaj"
$ws24.Range("P3").Value = 3

$ws24.Activate()
$ws24.Range("A3").Select()
